# "Fixed missing items due to empty line"
#
# The sheet "VINI, DISTILLATI & CO." had a stray empty row (row 24) sitting
# in the middle of the drink list, right before "TERREBLU TORBATA" (the
# first Grappa entry). Because of that blank row, the list looked like it
# was missing an item. The fix removes that empty row so every item shifts
# up by one row.

$wb = $excel.ActiveWorkbook

# Leave the previously-active sheet's selection as it will end up once it
# is no longer the active tab.
$ws1 = $wb.Worksheets.Item("LE NOSTRE BIRRE")
$ws1.Range("E23").Select()

# Remove the stray empty row from the drinks list and make this sheet the
# active one (mirrors the tabSelected/activeTab change in the workbook).
$ws = $wb.Worksheets.Item("VINI, DISTILLATI & CO.")
$ws.Rows.Item(24).Delete()
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("A24").Select()
